$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 ("mf" -> "NN") with new metric values
$ws.Range("A2").Value = "NN"
$ws.Range("B2").Value = 4.431730128664218
$ws.Range("C2").Value = 3.527280469920079
$ws.Range("D2").Value = 0.06518861681005957
$ws.Range("E2").Value = 0.2813446562071336
$ws.Range("F2").Value = 0.9013412816691506
$ws.Range("G2").Value = 0.1679717659980404
$ws.Range("H2").Value = 0.2482472170010725
$ws.Range("I2").Value = 0.9046274984948548
$ws.Range("J2").Value = 0.1592793133421149
$ws.Range("K2").Value = 0.9526080476900149
$ws.Range("L2").Value = 0.6016746531844853
$ws.Range("M2").Value = 0.06518861681005957

# Remove row 3 (the "mmr" row) entirely
$ws.Rows.Item(3).Delete()
